# Auto-generated script applying market-data value updates per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 116.5
$ws.Range("J9").Value = 46.666668
$ws.Range("L9").Value = 46.666668
$ws.Range("N9").Value = -384.666668
$ws.Range("H33").Value = 450.44446
$ws.Range("I33").Value = 521.0714
$ws.Range("J33").Value = 203.25
$ws.Range("K33").Value = 521.0714
$ws.Range("L33").Value = 203.25
$ws.Range("M33").Value = -292.0714
$ws.Range("N33").Value = -661.25
$ws.Range("H40").Value = 3457.3845
$ws.Range("J40").Value = 6644.25
$ws.Range("L40").Value = 6644.25
$ws.Range("N40").Value = -6994.25
$ws.Range("H113").Value = 5806.4614
$ws.Range("I113").Value = 2676.8
$ws.Range("K113").Value = 2676.8
$ws.Range("M113").Value = 577.1999999999998
$ws.Range("H132").Value = 15425.892
$ws.Range("I132").Value = 1272.72
$ws.Range("K132").Value = 3818.16
$ws.Range("M132").Value = -1288.16

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10056.459
$ws.Range("I32").Value = 9096.727999999999
$ws.Range("K32").Value = 9096.727999999999
$ws.Range("M32").Value = -8809.727999999999
$ws.Range("H45").Value = 2780.2
$ws.Range("I45").Value = 2227.6
$ws.Range("K45").Value = 2227.6
$ws.Range("M45").Value = -1850.6
$ws.Range("H102").Value = 4083.3333
$ws.Range("I102").Value = 4666.6665
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 4666.6665
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -3044.6665
$ws.Range("N102").Value = -6744
$ws.Range("H122").Value = 3060.1
$ws.Range("I122").Value = 1979.1875
$ws.Range("K122").Value = 5937.5625
$ws.Range("M122").Value = -3487.5625
$ws.Range("H132").Value = 1322.6666
$ws.Range("I132").Value = 984.25
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 2952.75
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -422.75
$ws.Range("N132").Value = -11058.5
$ws.Range("H139").Value = 77666.664
$ws.Range("J139").Value = 76500
$ws.Range("L139").Value = 76500
$ws.Range("N139").Value = -86780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3284.4644
$ws.Range("J20").Value = 5970
$ws.Range("L20").Value = 5970
$ws.Range("N20").Value = -6464
$ws.Range("H86").Value = 3958.818
$ws.Range("I86").Value = 1799.1875
$ws.Range("J86").Value = 9717.833000000001
$ws.Range("K86").Value = 1799.1875
$ws.Range("L86").Value = 9717.833000000001
$ws.Range("M86").Value = -676.1875
$ws.Range("N86").Value = -11963.833
$ws.Range("H89").Value = 3958.818
$ws.Range("I89").Value = 1799.1875
$ws.Range("J89").Value = 9717.833000000001
$ws.Range("K89").Value = 8995.9375
$ws.Range("L89").Value = 48589.165
$ws.Range("M89").Value = -3379.9375
$ws.Range("N89").Value = -59821.165
$ws.Range("H94").Value = 3681
$ws.Range("I94").Value = 3631
$ws.Range("J94").Value = 3737.25
$ws.Range("K94").Value = 3631
$ws.Range("L94").Value = 3737.25
$ws.Range("M94").Value = -3180
$ws.Range("N94").Value = -4639.25
$ws.Range("H105").Value = 2633.5
$ws.Range("I105").Value = 3800.3333
$ws.Range("J105").Value = 1466.6666
$ws.Range("K105").Value = 3800.3333
$ws.Range("L105").Value = 1466.6666
$ws.Range("M105").Value = -2053.3333
$ws.Range("N105").Value = -4960.6666
$ws.Range("H134").Value = 3282.4187
$ws.Range("I134").Value = 2004.697
$ws.Range("K134").Value = 6014.090999999999
$ws.Range("M134").Value = -3479.090999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 781.36365
$ws.Range("I16").Value = 602.8333
$ws.Range("K16").Value = 602.8333
$ws.Range("M16").Value = -315.8333
$ws.Range("H86").Value = 35190.145
$ws.Range("I86").Value = 45621.125
$ws.Range("J86").Value = 21282.166
$ws.Range("K86").Value = 45621.125
$ws.Range("L86").Value = 21282.166
$ws.Range("M86").Value = -44498.125
$ws.Range("N86").Value = -23528.166
$ws.Range("H89").Value = 35190.145
$ws.Range("I89").Value = 45621.125
$ws.Range("J89").Value = 21282.166
$ws.Range("K89").Value = 228105.625
$ws.Range("L89").Value = 106410.83
$ws.Range("M89").Value = -222489.625
$ws.Range("N89").Value = -117642.83
$ws.Range("H99").Value = 10689680
$ws.Range("I99").Value = 1877866.4
$ws.Range("K99").Value = 1877866.4
$ws.Range("M99").Value = -1876368.4
$ws.Range("H113").Value = 781.36365
$ws.Range("I113").Value = 602.8333
$ws.Range("K113").Value = 602.8333
$ws.Range("M113").Value = 1567.1667
$ws.Range("H122").Value = 322623.62
$ws.Range("I122").Value = 730519.8
$ws.Range("J122").Value = 5371.0557
$ws.Range("K122").Value = 2191559.4
$ws.Range("L122").Value = 16113.1671
$ws.Range("M122").Value = -2189109.4
$ws.Range("N122").Value = -21013.1671
$ws.Range("H126").Value = 10689680
$ws.Range("I126").Value = 1877866.4
$ws.Range("K126").Value = 5633599.199999999
$ws.Range("M126").Value = -5631129.199999999
$ws.Range("H134").Value = 3865.762
$ws.Range("I134").Value = 2373.111
$ws.Range("K134").Value = 7119.333
$ws.Range("M134").Value = -4584.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 918.4
$ws.Range("I7").Value = 1124.4667
$ws.Range("J7").Value = 300.2
$ws.Range("K7").Value = 3373.4001
$ws.Range("L7").Value = 900.5999999999999
$ws.Range("M7").Value = -3261.4001
$ws.Range("N7").Value = -1124.6
$ws.Range("H16").Value = 336.66666
$ws.Range("I16").Value = 336.66666
$ws.Range("K16").Value = 1009.99998
$ws.Range("M16").Value = -836.9999799999999
$ws.Range("H17").Value = 1100
$ws.Range("I17").Value = 350
$ws.Range("K17").Value = 1050
$ws.Range("M17").Value = -881
$ws.Range("H20").Value = 767
$ws.Range("I20").Value = 901
$ws.Range("J20").Value = 700
$ws.Range("K20").Value = 2703
$ws.Range("L20").Value = 2100
$ws.Range("M20").Value = -2476
$ws.Range("N20").Value = -2554
$ws.Range("H24").Value = 169
$ws.Range("I24").Value = 169
$ws.Range("K24").Value = 507
$ws.Range("M24").Value = -277
$ws.Range("H121").Value = 1985.4
$ws.Range("J121").Value = 2138.2222
$ws.Range("L121").Value = 6414.6666
$ws.Range("N121").Value = -9034.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58365.15
$ws.Range("I80").Value = 87869.16
$ws.Range("K80").Value = 87869.16
$ws.Range("M80").Value = -86871.16
$ws.Range("H83").Value = 58365.15
$ws.Range("I83").Value = 87869.16
$ws.Range("K83").Value = 439345.8
$ws.Range("M83").Value = -434353.8
$ws.Range("H102").Value = 1011.4286
$ws.Range("I102").Value = 897.9167
$ws.Range("J102").Value = 1259.091
$ws.Range("K102").Value = 897.9167
$ws.Range("L102").Value = 1259.091
$ws.Range("M102").Value = 724.0833
$ws.Range("N102").Value = -4503.091
$ws.Range("H122").Value = 5050.0454
$ws.Range("I122").Value = 1465.6364
$ws.Range("J122").Value = 8634.454
$ws.Range("K122").Value = 4396.9092
$ws.Range("L122").Value = 25903.362
$ws.Range("M122").Value = -1946.9092
$ws.Range("N122").Value = -30803.362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 319.77777
$ws.Range("I16").Value = 319.77777
$ws.Range("K16").Value = 319.77777
$ws.Range("M16").Value = -149.77777
$ws.Range("H22").Value = 563.5
$ws.Range("I22").Value = 576.125
$ws.Range("J22").Value = 546.6667
$ws.Range("K22").Value = 576.125
$ws.Range("L22").Value = 546.6667
$ws.Range("M22").Value = -281.125
$ws.Range("N22").Value = -1136.6667
$ws.Range("H27").Value = 563.5
$ws.Range("I27").Value = 576.125
$ws.Range("J27").Value = 546.6667
$ws.Range("K27").Value = 576.125
$ws.Range("L27").Value = 546.6667
$ws.Range("M27").Value = -469.125
$ws.Range("N27").Value = -760.6667
$ws.Range("H40").Value = 11126.177
$ws.Range("I40").Value = 18989.166
$ws.Range("J40").Value = 6837.273
$ws.Range("K40").Value = 18989.166
$ws.Range("L40").Value = 6837.273
$ws.Range("M40").Value = -18853.166
$ws.Range("N40").Value = -7109.273
$ws.Range("H46").Value = 935.3333
$ws.Range("I46").Value = 895.1429000000001
$ws.Range("J46").Value = 1076
$ws.Range("K46").Value = 895.1429000000001
$ws.Range("L46").Value = 1076
$ws.Range("M46").Value = -707.1429000000001
$ws.Range("N46").Value = -1452
$ws.Range("H55").Value = 1144.7273
$ws.Range("I55").Value = 1199.2106
$ws.Range("K55").Value = 1199.2106
$ws.Range("M55").Value = -1026.2106
$ws.Range("H82").Value = 4760.65
$ws.Range("J82").Value = 7300
$ws.Range("L82").Value = 7300
$ws.Range("N82").Value = -8022
$ws.Range("H85").Value = 4760.65
$ws.Range("J85").Value = 7300
$ws.Range("L85").Value = 7300
$ws.Range("N85").Value = -9796
$ws.Range("H93").Value = 5552.6
$ws.Range("J93").Value = 1892
$ws.Range("L93").Value = 1892
$ws.Range("N93").Value = -4388
$ws.Range("H122").Value = 5994.1816
$ws.Range("I122").Value = 2804
$ws.Range("K122").Value = 8412
$ws.Range("M122").Value = -5962

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2421.7144
$ws.Range("I62").Value = 2075.5
$ws.Range("J62").Value = 2883.3333
$ws.Range("K62").Value = 2075.5
$ws.Range("L62").Value = 2883.3333
$ws.Range("M62").Value = -1451.5
$ws.Range("N62").Value = -4131.3333
$ws.Range("H65").Value = 2421.7144
$ws.Range("I65").Value = 2075.5
$ws.Range("J65").Value = 2883.3333
$ws.Range("K65").Value = 10377.5
$ws.Range("L65").Value = 14416.6665
$ws.Range("M65").Value = -7257.5
$ws.Range("N65").Value = -20656.6665
$ws.Range("H96").Value = 1997
$ws.Range("J96").Value = 2100
$ws.Range("L96").Value = 2100
$ws.Range("N96").Value = -4846
$ws.Range("H132").Value = 1676.6
$ws.Range("I132").Value = 1323
$ws.Range("J132").Value = 3975
$ws.Range("K132").Value = 3969
$ws.Range("L132").Value = 11925
$ws.Range("M132").Value = -1439
$ws.Range("N132").Value = -16985
